# Auto-generated Excel COM-interop script
# Updates numeric market-data cells (H,I,J,K,L,M,N) across all 8 sheets
# to match the refreshed values pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 946.1539
$ws.Range("I18").Value = 946.1539
$ws.Range("K18").Value = 946.1539
$ws.Range("M18").Value = -662.1539
$ws.Range("H33").Value = 461.4
$ws.Range("I33").Value = 461.4
$ws.Range("K33").Value = 461.4
$ws.Range("M33").Value = -232.4
$ws.Range("H138").Value = 2330254.2
$ws.Range("J138").Value = 6900.6045
$ws.Range("L138").Value = 20701.8135
$ws.Range("N138").Value = -30981.8135

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 130.5
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 111
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 111
$ws.Range("M4").Value = -34
$ws.Range("N4").Value = -343
$ws.Range("H32").Value = 17979.55
$ws.Range("I32").Value = 13421.741
$ws.Range("J32").Value = 58999.832
$ws.Range("K32").Value = 13421.741
$ws.Range("L32").Value = 58999.832
$ws.Range("M32").Value = -13134.741
$ws.Range("N32").Value = -59573.832
$ws.Range("H102").Value = 3600
$ws.Range("I102").Value = 1466.6666
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 1466.6666
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 155.3334
$ws.Range("N102").Value = -13244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 57781.477
$ws.Range("I75").Value = 8975.666999999999
$ws.Range("J75").Value = 111024.18
$ws.Range("K75").Value = 8975.666999999999
$ws.Range("L75").Value = 111024.18
$ws.Range("M75").Value = -8039.666999999999
$ws.Range("N75").Value = -112896.18
$ws.Range("H78").Value = 57781.477
$ws.Range("I78").Value = 8975.666999999999
$ws.Range("J78").Value = 111024.18
$ws.Range("K78").Value = 26927.001
$ws.Range("L78").Value = 333072.54
$ws.Range("M78").Value = -22247.001
$ws.Range("N78").Value = -342432.54
$ws.Range("H94").Value = 1536.6666
$ws.Range("I94").Value = 1393.3334
$ws.Range("J94").Value = 1608.3334
$ws.Range("K94").Value = 1393.3334
$ws.Range("L94").Value = 1608.3334
$ws.Range("M94").Value = -942.3334
$ws.Range("N94").Value = -2510.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 158.5
$ws.Range("I7").Value = 83.666664
$ws.Range("K7").Value = 83.666664
$ws.Range("M7").Value = 29.333336
$ws.Range("H11").Value = 25200
$ws.Range("I11").Value = 25200
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 25200
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -25060
$ws.Range("N11").ClearContents()
$ws.Range("H31").Value = 3632.468
$ws.Range("I31").Value = 2292.6667
$ws.Range("J31").Value = 5030.522
$ws.Range("K31").Value = 2292.6667
$ws.Range("L31").Value = 5030.522
$ws.Range("M31").Value = -1997.6667
$ws.Range("N31").Value = -5620.522
$ws.Range("H34").Value = 3632.468
$ws.Range("I34").Value = 2292.6667
$ws.Range("J34").Value = 5030.522
$ws.Range("K34").Value = 2292.6667
$ws.Range("L34").Value = 5030.522
$ws.Range("M34").Value = -2090.6667
$ws.Range("N34").Value = -5434.522
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 23000
$ws.Range("H58").Value = 1371.9642
$ws.Range("I58").Value = 1404.7916
$ws.Range("J58").Value = 1175
$ws.Range("K58").Value = 1404.7916
$ws.Range("L58").Value = 1175
$ws.Range("M58").Value = -1201.7916
$ws.Range("N58").Value = -1581
$ws.Range("H60").Value = 15000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 23000
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488
$ws.Range("H74").Value = 34000
$ws.Range("J74").Value = 34000
$ws.Range("L74").Value = 34000
$ws.Range("N74").Value = -35748
$ws.Range("H77").Value = 34000
$ws.Range("J77").Value = 34000
$ws.Range("L77").Value = 102000
$ws.Range("N77").Value = -110736
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 1534.2727
$ws.Range("I105").Value = 1011
$ws.Range("K105").Value = 1011
$ws.Range("M105").Value = 736
$ws.Range("H132").Value = 2156.125
$ws.Range("I132").Value = 2156.125
$ws.Range("K132").Value = 6468.375
$ws.Range("M132").Value = -3938.375
$ws.Range("H136").Value = 1371.9642
$ws.Range("I136").Value = 1404.7916
$ws.Range("J136").Value = 1175
$ws.Range("K136").Value = 4214.3748
$ws.Range("L136").Value = 3525
$ws.Range("M136").Value = -1664.3748
$ws.Range("N136").Value = -8625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2382838
$ws.Range("I129").Value = 733.3333
$ws.Range("J129").Value = 2779855.5
$ws.Range("K129").Value = 2199.9999
$ws.Range("L129").Value = 8339566.5
$ws.Range("M129").Value = 2800.0001
$ws.Range("N129").Value = -8349566.5
$ws.Range("H131").Value = 817.55554
$ws.Range("J131").Value = 859.6445
$ws.Range("L131").Value = 2578.9335
$ws.Range("N131").Value = -12658.9335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 55000000
$ws.Range("I11").Value = 55000000
$ws.Range("K11").Value = 55000000
$ws.Range("M11").Value = -54999861
$ws.Range("H124").Value = 69780
$ws.Range("J124").Value = 69780
$ws.Range("L124").Value = 69780
$ws.Range("N124").Value = -79600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4888.8887
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 4800
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 4800
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -5024
$ws.Range("H102").Value = 49561
$ws.Range("J102").Value = 49561
$ws.Range("L102").Value = 49561
$ws.Range("N102").Value = -56051
$ws.Range("H111").Value = 48385.332
$ws.Range("J111").Value = 48385.332
$ws.Range("L111").Value = 48385.332
$ws.Range("N111").Value = -56565.332
$ws.Range("H122").Value = 1182178.6
$ws.Range("I122").Value = 1639692.4
$ws.Range("J122").Value = 5714.857
$ws.Range("K122").Value = 4919077.199999999
$ws.Range("L122").Value = 17144.571
$ws.Range("M122").Value = -4916627.199999999
$ws.Range("N122").Value = -22044.571
$ws.Range("H126").Value = 4888.8887
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 4800
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 14400
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -19340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 24114
$ws.Range("J64").Value = 24114
$ws.Range("L64").Value = 24114
$ws.Range("N64").Value = -24610
$ws.Range("H67").Value = 24114
$ws.Range("J67").Value = 24114
$ws.Range("L67").Value = 24114
$ws.Range("N67").Value = -25830
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360
$ws.Range("H81").Value = 252749.75
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 252749.75
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608
$ws.Range("H102").Value = 34000
$ws.Range("J102").Value = 34000
$ws.Range("L102").Value = 34000
$ws.Range("N102").Value = -40490
$ws.Range("H109").Value = 44251.332
$ws.Range("J109").Value = 44251.332
$ws.Range("L109").Value = 44251.332
$ws.Range("N109").Value = -47025.332
